# Update workbook to add data "through October 19" (commit: "Add data for 2022-10-27")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (Through 2022-10-18 -> Through 2022-10-19)
$ws.Name = "Through 2022-10-19"

# 2. Update the column header text for the "October 2022" column (B1)
$ws.Range("B1").Value = "October 2022 (through October 19)"

# 3. Update / add the individual neighborhood counts
$ws.Range("L2").Value = 17     # Garfield Park, October 2021: 13 -> 17
$ws.Range("V2").Value = 12     # Garfield Park, October 2020: 11 -> 12
$ws.Range("B3").Value = 4      # Humboldt Park, October 2022: 3 -> 4
$ws.Range("L4").Value = 6      # South Shore, October 2021: 5 -> 6
$ws.Range("AZ7").Value = 3     # Englewood, October 2017: 2 -> 3
$ws.Range("V8").Value = 5      # Washington Heights, October 2020: 4 -> 5
$ws.Range("B10").Value = 4     # New City, October 2022: 3 -> 4
$ws.Range("L16").Value = 1     # Little Italy, UIC, October 2021: (new) -> 1
$ws.Range("L23").Value = 3     # Grand Boulevard, October 2021: 2 -> 3
$ws.Range("B24").Value = 5     # Auburn Gresham, October 2022: 4 -> 5
$ws.Range("L24").Value = 5     # Auburn Gresham, October 2021: 4 -> 5
$ws.Range("L36").Value = 1     # Washington Park, October 2021: (new) -> 1
$ws.Range("L42").Value = 2     # Avondale, October 2021: 1 -> 2
$ws.Range("V43").Value = 3     # Bridgeport, October 2020: 2 -> 3
$ws.Range("BT43").Value = 1    # Bridgeport, October 2015: (new) -> 1
$ws.Range("L47").Value = 2     # Chinatown, October 2021: 1 -> 2
$ws.Range("L49").Value = 1     # Near South Side, October 2021: (new) -> 1
$ws.Range("V70").Value = 2     # Gage Park, October 2020: 1 -> 2
$ws.Range("AP95").Value = 1    # United Center, October 2018: (new) -> 1
$ws.Range("BJ98").Value = 1    # Woodlawn, October 2016: (new) -> 1
